$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: new execution timestamp/trace/workflow id, date moved to 2023-04-04 (45020)
$ws.Range("F8").Value = 45020
$ws.Range("G8").Value = "2023-04-04T08:08:12Z"
$ws.Range("H8").Value = "10c0f47781b3d506"
$ws.Range("I8").Value = "2.16.840.1.113883.2.9.2.30.fb45f3577b696e54209bc2b936b93d8a70a7fa3458a50e694282b1d5799afd9f.d1e963401f^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"

# Row 9
$ws.Range("F9").Value = 45020
$ws.Range("G9").Value = "2023-04-04T08:10:56Z"
$ws.Range("H9").Value = "889bd3999b120994"
$ws.Range("I9").Value = "2.16.840.1.113883.2.9.2.30.64f74863595b6c2e7c8715ebaf432ceabb5091659c9369541a2641e60a80aeb3.4e77841d9f^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"

# Row 10
$ws.Range("F10").Value = 45020
$ws.Range("G10").Value = "2023-04-04T08:13:05Z"
$ws.Range("H10").Value = "1805a5034fd5b06d"
$ws.Range("I10").Value = "2.16.840.1.113883.2.9.2.30.2d257b5b8c756775146115c50b0d97851c0a657dfb3d22fe8561cfe9c1f1e2b4.babb7c65c2^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"

# Row 11
$ws.Range("F11").Value = 45020
$ws.Range("G11").Value = "2023-04-04T08:18:58Z"
$ws.Range("H11").Value = "9d2a04ab2ca6ec02"
$ws.Range("I11").Value = "2.16.840.1.113883.2.9.2.30.b3f88f87fbd4d91f876411f8dd396e4d8491a8506178c07d2b78c911a0f8b320.94f09de098^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"

# Row 12 - note: TRACEID and WORKFLOWINSTANCEID are entered before the TIMESTAMP for this row
$ws.Range("F12").Value = 45020
$ws.Range("H12").Value = "07ec4b2c686c1c09"
$ws.Range("I12").Value = "2.16.840.1.113883.2.9.2.30.2b51ae544a471226d473780a5f9f0ac3714f12b3c9f8a93d4efe1d102774d3b3.362a84b5eb^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("G12").Value = "22023-04-04T08:21:07Z"

# Update the view selection to match the corrected/reviewed cell
$ws.Range("I13").Select() | Out-Null
